$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full updated data table (rows 2-7), columns A-T
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E Ligand-expressing cells ... T Edge total expression derived specificity
$data = @(
    @("ECs","Hsp90aa1","Cftr","ECs",3,1,155.2826053333333,465.847816,0.4511422706622114,0.4511422706622115,1,0.3333333333333333,0.02245766666666667,0.067373,0.1039642614653416,0.1039642614653416,3.487284989707555,31.385564907368,0.04690267298519407,0.04690267298519408),
    @("ECs","Hsp90aa1","Cftr","sCs",3,1,155.2826053333333,465.847816,0.4511422706622114,0.4511422706622115,3,1,0.1935556666666667,0.580667,0.8960357385346585,0.8960357385346583,30.05582819703022,270.502453773272,0.4042395976770173,0.4042395976770174),
    @("FAPs","Hsp90aa1","Cftr","ECs",3,1,117.8650183333333,353.595055,0.3424330232507294,0.3424330232507295,1,0.3333333333333333,0.02245766666666667,0.067373,0.1039642614653416,0.1039642614653416,2.646973293390556,23.822759640515,0.03560079636360625,0.03560079636360625),
    @("FAPs","Hsp90aa1","Cftr","sCs",3,1,117.8650183333333,353.595055,0.3424330232507294,0.3424330232507295,3,1,0.1935556666666667,0.580667,0.8960357385346585,0.8960357385346583,22.81344220018723,205.320979801685,0.3068322268871232,0.3068322268871232),
    @("sCs","Hsp90aa1","Cftr","ECs",3,1,71.05112566666666,213.153377,0.2064247060870591,0.2064247060870591,1,0.3333333333333333,0.02245766666666667,0.067373,0.1039642614653416,0.1039642614653416,1.595642496513444,14.360782468621,0.02146079211654131,0.02146079211654131),
    @("sCs","Hsp90aa1","Cftr","sCs",3,1,71.05112566666666,213.153377,0.2064247060870591,0.2064247060870591,3,1,0.1935556666666667,0.580667,0.8960357385346585,0.8960357385346583,13.75234799582878,123.771131962459,0.1849639139705178,0.1849639139705178)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
